# Refresh crypto price/volume table (coinranking.com feed) -- Aug 31 2024 run
# Two pairs of adjacent rows also swapped rank order (30/31 and 38/39, 43/44).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '58.871.02'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.32%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.496.55'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.32%  '

# Row 4
$ws.Range('E4').Value = '  -0.21%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '532.78'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.04%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '134.89'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.03%  '

# Row 7
$ws.Range('E7').Value = '  +0.08%  '

# Row 8
$ws.Range('E8').Value = '  +1.11%  '

# Row 9
$ws.Range('E9').Value = '  +1.00%  '

# Row 10
$ws.Range('E10').Value = '  -1.03%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.37'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.51%  '

# Row 12
$ws.Range('E12').Value = '  +0.42%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.938.45'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -0.03%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '58.801.24'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.17%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '22.68'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -1.46%  '

# Row 16
$ws.Range('E16').Value = '  -0.55%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.490.67'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.58%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.03'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.41%  '

# Row 19
$ws.Range('E19').Value = '  +0.13%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '322.55'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.40%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.00'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.04%  '

# Row 22
$ws.Range('E22').Value = '  +1.35%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.98'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +2.51%  '

# Row 24
$ws.Range('E24').Value = '  +0.64%  '

# Row 25
$ws.Range('E25').Value = '  -0.37%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.01'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.84%  '

# Row 27
$ws.Range('E27').Value = '  -0.51%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0₃0757'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.38%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '170.08'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.15%  '

# Row 30
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.74'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.22%  '

# Row 31
$ws.Range('B31').Value = 'Aptos'
$ws.Range('C31').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.43'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -4.23%  '

# Row 32
$ws.Range('E32').Value = '  +1.24%  '

# Row 33
$ws.Range('E33').Value = '  +0.02%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.31'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.60%  '

# Row 35
$ws.Range('E35').Value = '  -2.03%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '4.03'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.09%  '

# Row 37
$ws.Range('E37').Value = '  -2.68%  '

# Row 38
$ws.Range('B38').Value = 'Filecoin'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.56'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.16%  '

# Row 39
$ws.Range('B39').Value = 'SuiNetwork'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.798'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -1.98%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '280.84'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.25%  '

# Row 41
$ws.Range('E41').Value = '  +0.34%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '4.98'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -4.96%  '

# Row 43
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '10.91'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.62%  '

# Row 44
$ws.Range('B44').Value = 'Aave'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '129.43'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.28%  '

# Row 45
$ws.Range('E45').Value = '  +0.18%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0923'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +0.10%  '

# Row 47
$ws.Range('E47').Value = '  -2.29%  '

# Row 48
$ws.Range('E48').Value = '  -1.32%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '17.21'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.61%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.751.17'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.76%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.982'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -0.48%  '

